$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("layoffs_by_IPOstatus")
$ws.Activate()

$ws.Range("C2").Value = "num_layoffs"
$ws.Range("C3").Value = 24758
$ws.Range("C4").Value = 60051

$ws.Range("C2").Font.Bold = $true

$ws.Range("D11").Select()
